$d = $word.ActiveDocument

# --- 1. Merge the "EXP NO: 0" + "8" runs into a single "EXP NO: 08" run ---
# The visible text is already "EXP NO: 08" (split across two runs with identical
# bold formatting); running a Find/Replace over it rewrites it as one run.
$d.Content.Find.Execute("EXP NO: 08", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "EXP NO: 08", 2) | Out-Null

# --- 2. Materialize even / default / first page headers & footers ---
# Touching Section.Headers/Footers ranges is what mints header1-3.xml,
# footer1-3.xml, footnotes.xml and endnotes.xml and wires up the six
# header/footerReference entries in sectPr.
$sec = $d.Sections(1)

# Leave the "even" and "first" footers/headers blank.
$sec.Footers(3).Range.Text = ""   # wdHeaderFooterEvenPages  -> footer1.xml
$sec.Footers(2).Range.Text = ""   # wdHeaderFooterFirstPage  -> footer3.xml
$sec.Headers(1).Range.Text = ""   # wdHeaderFooterPrimary    -> header2.xml
$sec.Headers(2).Range.Text = ""   # wdHeaderFooterFirstPage  -> header3.xml
$sec.Headers(3).Range.Text = ""   # wdHeaderFooterEvenPages  -> header1.xml

# Default (primary) footer gets the new signature text.
$defFooter = $sec.Footers(1)      # wdHeaderFooterPrimary -> footer2.xml
$defFooter.Range.Text = "192110493                                                                                                                  KEERTHI"
$defFooter.Range.LanguageID = "en-IN"
